$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E -> B:F.
$ws.Columns("A").Insert()

# New header cell for the inserted ID column; match the header formatting
# (bold/bordered/centered style) already used by the rest of row 1.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Row labels (specimen IDs) for the newly inserted column A.
$ids = @(
    @(2,  "Hb 2"),
    @(3,  "Hb 3"),
    @(4,  "S 24"),
    @(5,  "S 28"),
    @(6,  "Hb 107"),
    @(7,  "Hb 66"),
    @(8,  "Hb 69"),
    @(9,  "Hb 95"),
    @(10, "Hb 99"),
    @(11, "Hb 92"),
    @(12, "Hb 40"),
    @(13, "Hb 41"),
    @(14, "S 11"),
    @(15, "Hb 57"),
    @(16, "S 21"),
    @(17, "S 22"),
    @(18, "S 3"),
    @(19, "S 4"),
    @(20, "S 5"),
    @(21, "Hb 74"),
    @(22, "Hb 79"),
    @(23, "Hb 32"),
    @(24, "S 15"),
    @(25, "S 16")
)

foreach ($pair in $ids) {
    $row = $pair[0]
    $label = $pair[1]
    $ws.Cells.Item($row, 1).Value = $label
}
